# MAJ DIAPO REVUE 1
# Repositions several shapes on slide 6 and appends two bullet lines to the
# "Rectangle 32" task list (also resized/repositioned).
#
# NOTE: this runtime stores Shape.Left/Top/Width/Height internally as
# single-precision (float32) point values and truncates (rather than
# rounds) when converting back to EMU on save. Adding half an EMU's worth
# of points before assigning compensates for that truncation so the saved
# XML round-trips to the exact target EMU value.
$EmuPerPoint = 12700

function ConvertTo-Pt($emu) {
    return ([double]$emu + 0.5) / $EmuPerPoint
}

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Set-ShapePos($slide, $id, $x, $y) {
    $sh = Get-ShapeById $slide $id
    $sh.Left = ConvertTo-Pt $x
    $sh.Top = ConvertTo-Pt $y
}

function Set-ShapePosExt($slide, $id, $x, $y, $cx, $cy) {
    $sh = Get-ShapeById $slide $id
    $sh.Left = ConvertTo-Pt $x
    $sh.Top = ConvertTo-Pt $y
    $sh.Width = ConvertTo-Pt $cx
    $sh.Height = ConvertTo-Pt $cy
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Rectangle 1 (id=2)
Set-ShapePos $s 2 143436 1886284
# Graphique 4 (id=5)
Set-ShapePos $s 5 143436 1011919
# Graphique 6 (id=7)
Set-ShapePos $s 7 361891 3420089
# Graphique 8 (id=9)
Set-ShapePos $s 9 388803 3869080
# Graphique 10 (id=11)
Set-ShapePos $s 11 453248 2258368
# Graphique 12 (id=13)
Set-ShapePos $s 13 244158 1265935
# Graphique 14 (id=15)
Set-ShapePos $s 15 361891 3033807
# Graphique 16 (id=17)
Set-ShapePos $s 17 426337 2512335
# Graphique 18 (id=19)
Set-ShapePos $s 19 387359 2730593
# ZoneTexte 19 (id=20)
Set-ShapePos $s 20 1101829 1279589
# Rectangle 21 (id=22)
Set-ShapePos $s 22 5137921 1886283
# Connecteur droit avec flèche 22 (id=23)
Set-ShapePos $s 23 3344249 3182475
# ZoneTexte 27 (id=28)
Set-ShapePos $s 28 7261797 1358501
# Picture 2 (id=27)
Set-ShapePos $s 27 5137921 1412000
# Rectangle 30 (id=31)
Set-ShapePos $s 31 5195962 1476314

# Rectangle 32 (id=33) - moves AND grows
Set-ShapePosExt $s 33 1701832 5148140 5928687 1612174

# ZoneTexte 33 (id=34)
Set-ShapePos $s 34 2177300 4677740
# Graphique 34 (id=35)
Set-ShapePos $s 35 1701832 4586095

# Append two new bullet paragraphs to Rectangle 32's text, matching the
# style (bullet / indent / run formatting) of the preceding paragraph.
$rect32 = Get-ShapeById $s 33
$tr = $rect32.TextFrame.TextRange
$tr.InsertAfter([char]13 + "Installation clé WIFI sur Raspberry") | Out-Null
$tr2 = $rect32.TextFrame.TextRange
$tr2.InsertAfter([char]13 + "Paramétrage de l’adaptateur RS232/RJ45") | Out-Null
